# Apply "Add data for 2021-12-20" update to the carjacking arrests workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2021-12-12"

# Row 10 (July) - 2021 column (T/U/V) updated
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 152
$ws.Range("V10").Value = 0.05

# Row 13 (November) - 2021 column (U/V) updated
$ws.Range("U13").Value = 195
$ws.Range("V13").Value = 0.025

# Row 14 (December, through date changes) - label + all year columns updated
$ws.Range("A14").Value = "December (through 12-12)"

$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 0.1667

$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 30
$ws.Range("G14").Value = 0.1176

$ws.Range("I14").Value = 37
$ws.Range("J14").Value = 0.075

$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 24
$ws.Range("M14").Value = 0.1111

$ws.Range("O14").Value = 15
$ws.Range("P14").Value = 0.1667

$ws.Range("R14").Value = 57
$ws.Range("S14").Value = 0.05

$ws.Range("U14").Value = 90

# Row 15 (Total) - recomputed totals across all years
$ws.Range("C15").Value = 268
$ws.Range("D15").Value = 0.1155

$ws.Range("E15").Value = 64
$ws.Range("F15").Value = 533
$ws.Range("G15").Value = 0.1072

$ws.Range("I15").Value = 795

$ws.Range("K15").Value = 77
$ws.Range("L15").Value = 632
$ws.Range("M15").Value = 0.1086

$ws.Range("O15").Value = 495
$ws.Range("P15").Value = 0.1033

$ws.Range("R15").Value = 1257
$ws.Range("S15").Value = 0.0506

$ws.Range("T15").Value = 100
$ws.Range("U15").Value = 1634
$ws.Range("V15").Value = 0.0577
